$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 922.4545
$ws.Range("I33").Value = 912.4
$ws.Range("K33").Value = 912.4
$ws.Range("M33").Value = -683.4
$ws.Range("H39").Value = 282.36365
$ws.Range("I39").Value = 290.7
$ws.Range("K39").Value = 872.0999999999999
$ws.Range("M39").Value = -576.0999999999999
$ws.Range("H41").Value = 780.1818
$ws.Range("I41").Value = 550.0
$ws.Range("J41").Value = 911.7143
$ws.Range("K41").Value = 550.0
$ws.Range("L41").Value = 911.7143
$ws.Range("M41").Value = -110.0
$ws.Range("N41").Value = -1791.7143
$ws.Range("H74").Value = 4145.579
$ws.Range("I74").Value = 4172.875
$ws.Range("K74").Value = 4172.875
$ws.Range("M74").Value = -3236.875
$ws.Range("H77").Value = 4145.579
$ws.Range("I77").Value = 4172.875
$ws.Range("K77").Value = 20864.375
$ws.Range("M77").Value = -16184.375
$ws.Range("H86").Value = 111113410.0
$ws.Range("J86").Value = 2399.25
$ws.Range("L86").Value = 2399.25
$ws.Range("N86").Value = -4645.25
$ws.Range("H89").Value = 111113410.0
$ws.Range("J89").Value = 2399.25
$ws.Range("L89").Value = 11996.25
$ws.Range("N89").Value = -23228.25
$ws.Range("H100").Value = 1859.0769
$ws.Range("I100").Value = 1459.5714
$ws.Range("K100").Value = 1459.5714
$ws.Range("M100").Value = -918.5714
$ws.Range("H116").Value = 7795.1943
$ws.Range("I116").Value = 7376.125
$ws.Range("J116").Value = 8130.45
$ws.Range("K116").Value = 7376.125
$ws.Range("L116").Value = 8130.45
$ws.Range("M116").Value = -3934.125
$ws.Range("N116").Value = -15014.45
$ws.Range("H132").Value = 4602.7847
$ws.Range("I132").Value = 2505.9167
$ws.Range("K132").Value = 7517.750100000001
$ws.Range("M132").Value = -4987.750100000001

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4995.035
$ws.Range("J32").Value = 4662.0713
$ws.Range("L32").Value = 4662.0713
$ws.Range("N32").Value = -5236.0713
$ws.Range("H61").Value = 4792.242
$ws.Range("J61").Value = 2964.3333
$ws.Range("L61").Value = 2964.3333
$ws.Range("N61").Value = -3388.3333
$ws.Range("H63").Value = 3392.1304
$ws.Range("I63").Value = 3271.7273
$ws.Range("J63").Value = 3502.5
$ws.Range("K63").Value = 3271.7273
$ws.Range("L63").Value = 3502.5
$ws.Range("M63").Value = -2585.7273
$ws.Range("N63").Value = -4874.5
$ws.Range("H64").Value = 49999.0
$ws.Range("J64").Value = 49999.0
$ws.Range("L64").Value = 49999.0
$ws.Range("N64").Value = -50495.0
$ws.Range("H66").Value = 3392.1304
$ws.Range("I66").Value = 3271.7273
$ws.Range("J66").Value = 3502.5
$ws.Range("K66").Value = 16358.6365
$ws.Range("L66").Value = 17512.5
$ws.Range("M66").Value = -12926.6365
$ws.Range("N66").Value = -24376.5
$ws.Range("H67").Value = 49999.0
$ws.Range("J67").Value = 49999.0
$ws.Range("L67").Value = 49999.0
$ws.Range("N67").Value = -51715.0
$ws.Range("H101").Value = 35000.0
$ws.Range("J101").Value = 35000.0
$ws.Range("L101").Value = 35000.0
$ws.Range("N101").Value = -41490.0
$ws.Range("H136").Value = 4792.242
$ws.Range("J136").Value = 2964.3333
$ws.Range("L136").Value = 8892.999899999999
$ws.Range("N136").Value = -13992.9999
$ws.Range("H62").Value = 0.0
$ws.Range("J62").Value = 0.0
$ws.Range("L62").Value = 0.0
$ws.Range("H65").Value = 0.0
$ws.Range("J65").Value = 0.0
$ws.Range("L65").Value = 0.0
$ws.Range("H110").Value = 4813.6
$ws.Range("I110").Value = 5689.75
$ws.Range("J110").Value = 1309.0
$ws.Range("K110").Value = 5689.75
$ws.Range("L110").Value = 1309.0
$ws.Range("M110").Value = -3644.75
$ws.Range("N110").Value = -5399.0
$ws.Range("N62").ClearContents()
$ws.Range("N65").ClearContents()

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 4905.385
$ws.Range("I80").Value = 13238.0
$ws.Range("J80").Value = 1202.0
$ws.Range("K80").Value = 13238.0
$ws.Range("L80").Value = 1202.0
$ws.Range("M80").Value = -12240.0
$ws.Range("N80").Value = -3198.0
$ws.Range("H83").Value = 4905.385
$ws.Range("I83").Value = 13238.0
$ws.Range("J83").Value = 1202.0
$ws.Range("K83").Value = 66190.0
$ws.Range("L83").Value = 6010.0
$ws.Range("M83").Value = -61198.0
$ws.Range("N83").Value = -15994.0
$ws.Range("H105").Value = 4316.9644
$ws.Range("I105").Value = 3020.5217
$ws.Range("K105").Value = 3020.5217
$ws.Range("M105").Value = -1273.5217
$ws.Range("H134").Value = 1933.2623
$ws.Range("I134").Value = 1950.5
$ws.Range("K134").Value = 5851.5
$ws.Range("M134").Value = -3316.5

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 5330.25
$ws.Range("J94").Value = 4950.2856
$ws.Range("L94").Value = 4950.2856
$ws.Range("N94").Value = -5852.2856
$ws.Range("H105").Value = 2690.0557
$ws.Range("I105").Value = 1418.0
$ws.Range("K105").Value = 1418.0
$ws.Range("M105").Value = 329.0
$ws.Range("H107").Value = 1383.3
$ws.Range("J107").Value = 2192.6667
$ws.Range("L107").Value = 2192.6667
$ws.Range("N107").Value = -6032.6667
$ws.Range("H132").Value = 11101.277
$ws.Range("I132").Value = 3804.3635
$ws.Range("J132").Value = 22567.857
$ws.Range("K132").Value = 11413.0905
$ws.Range("L132").Value = 67703.571
$ws.Range("M132").Value = -8883.0905
$ws.Range("N132").Value = -72763.571
$ws.Range("H134").Value = 3088.0571
$ws.Range("I134").Value = 3451.3704
$ws.Range("J134").Value = 1861.875
$ws.Range("K134").Value = 10354.1112
$ws.Range("L134").Value = 5585.625
$ws.Range("M134").Value = -7819.111199999999
$ws.Range("N134").Value = -10655.625

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I2").Value = 64814860.0
$ws.Range("J2").Value = 84.4
$ws.Range("K2").Value = 64814860.0
$ws.Range("L2").Value = 84.4
$ws.Range("M2").Value = -64814747.0
$ws.Range("N2").Value = -310.4

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 889.38464
$ws.Range("I22").Value = 845.55554
$ws.Range("J22").Value = 988.0
$ws.Range("K22").Value = 845.55554
$ws.Range("L22").Value = 988.0
$ws.Range("M22").Value = -550.55554
$ws.Range("N22").Value = -1578.0
$ws.Range("H27").Value = 889.38464
$ws.Range("I27").Value = 845.55554
$ws.Range("J27").Value = 988.0
$ws.Range("K27").Value = 845.55554
$ws.Range("L27").Value = 988.0
$ws.Range("M27").Value = -738.55554
$ws.Range("N27").Value = -1202.0
$ws.Range("H40").Value = 6891.722
$ws.Range("I40").Value = 6876.2144
$ws.Range("J40").Value = 6946.0
$ws.Range("K40").Value = 6876.2144
$ws.Range("L40").Value = 6946.0
$ws.Range("M40").Value = -6740.2144
$ws.Range("N40").Value = -7218.0
$ws.Range("H100").Value = 62502800.0
$ws.Range("I100").Value = 66669452.0
$ws.Range("J100").Value = 3000.0
$ws.Range("K100").Value = 66669452.0
$ws.Range("L100").Value = 3000.0
$ws.Range("M100").Value = -66668911.0
$ws.Range("N100").Value = -4082.0

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1874.38
$ws.Range("I132").Value = 1883.75
$ws.Range("J132").Value = 1649.5
$ws.Range("K132").Value = 5651.25
$ws.Range("L132").Value = 4948.5
$ws.Range("M132").Value = -3121.25
$ws.Range("N132").Value = -10008.5
$ws.Range("H136").Value = 1665.6731
$ws.Range("I136").Value = 1686.34
$ws.Range("J136").Value = 1149.0
$ws.Range("K136").Value = 5059.02
$ws.Range("L136").Value = 3447.0
$ws.Range("M136").Value = -2509.02
$ws.Range("N136").Value = -8547.0

Write-Output "Applied all market-data updates across sheets"